# Auto-update draw results: append the 2025-11-04 Pick 3 row (row 49).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 49

# Date-like and purely-numeric strings must stay text (matches the rest of
# the sheet, which stores every column as text). A leading apostrophe forces
# Excel to keep the literal string instead of coercing it to a date serial /
# number; re-applying the "Normal" style afterwards drops the quote-prefix
# style Excel would otherwise tag the cell with, so the cell ends up with no
# explicit style - same as every other data cell on the sheet.
$ws.Cells.Item($row, 1).Value = "'2025-11-04"
$ws.Cells.Item($row, 1).Style = "Normal"

$ws.Cells.Item($row, 2).Value = "Pick 3"

$ws.Cells.Item($row, 3).Value = "'251104"
$ws.Cells.Item($row, 3).Style = "Normal"

$ws.Cells.Item($row, 4).Value = "5-5-5"

$ws.Cells.Item($row, 5).Value = "2025-11-04T21:39:26.107+04:00"
